$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.086.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.57%  "

$ws.Range("D3").Value = "'3.418.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.27%  "

$ws.Range("D5").Value = "'410.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Value = "'129.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("D7").Value = "'0.642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.85%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +4.64%  "

$ws.Range("E10").Value = "  +2.48%  "

$ws.Range("D11").Value = "'43.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("D12").Value = "'0.0000228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +39.15%  "

$ws.Range("D13").Value = "'9.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.35%  "

$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").Value = "'21.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.88%  "

$ws.Range("D16").Value = "'3.959.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.91%  "

$ws.Range("D17").Value = "'3.413.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.69%  "

$ws.Range("D18").Value = "'12.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.76%  "

$ws.Range("E19").Value = "  +6.88%  "

$ws.Range("D20").Value = "'62.042.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("D21").Value = "'502.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +59.22%  "

$ws.Range("D22").Value = "'93.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +10.74%  "

$ws.Range("D23").Value = "'3.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'13.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").Value = "'3.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.17%  "

$ws.Range("D26").Value = "'34.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.12%  "

$ws.Range("D27").Value = "'9.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.03%  "

$ws.Range("D28").Value = "'4.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("E29").Value = "  -3.44%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'12.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.87%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("E32").Value = "  -2.09%  "

$ws.Range("E33").Value = "  -0.99%  "

$ws.Range("D34").Value = "'42.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.70%  "

$ws.Range("E35").Value = "  +0.05%  "

$ws.Range("D36").Value = "'0.0511"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.68%  "

$ws.Range("D37").Value = "'54.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.81%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "

$ws.Range("D39").Value = "'0.138"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.51%  "

$ws.Range("E40").Value = "  +2.41%  "

$ws.Range("D41").Value = "'2.96"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("D42").Value = "'4.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.26%  "

$ws.Range("D43").Value = "'0.319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").Value = "'144.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("D45").Value = "'2.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +16.49%  "

$ws.Range("D46").Value = "'2.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.45%  "

$ws.Range("D47").Value = "'16.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.60%  "

$ws.Range("D48").Value = "'0.149"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.29%  "

$ws.Range("E49").Value = "  +5.36%  "

$ws.Range("D50").Value = "'113.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +34.46%  "

$ws.Range("E51").Value = "  +5.62%  "
